# Applies scheduled market-data refresh updates to the Gilgamesh Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 254.07143
$ws.Range("I9").Value = 354
$ws.Range("K9").Value = 354
$ws.Range("M9").Value = -185

$ws.Range("H19").Value = 1066.5
$ws.Range("J19").Value = 1637.2727
$ws.Range("L19").Value = 1637.2727
$ws.Range("N19").Value = -1987.2727

$ws.Range("H55").Value = 518.1429000000001
$ws.Range("J55").Value = 541.7857
$ws.Range("L55").Value = 541.7857
$ws.Range("N55").Value = -969.7857

$ws.Range("H116").Value = 9873.75
$ws.Range("I116").Value = 14750
$ws.Range("J116").Value = 4997.5
$ws.Range("K116").Value = 14750
$ws.Range("L116").Value = 4997.5
$ws.Range("M116").Value = -11308
$ws.Range("N116").Value = -11881.5

$ws.Range("H118").Value = 902.2
$ws.Range("I118").Value = 762.0833
$ws.Range("J118").Value = 1462.6666
$ws.Range("K118").Value = 2286.2499
$ws.Range("L118").Value = 4387.9998
$ws.Range("M118").Value = -629.2498999999998
$ws.Range("N118").Value = -7701.9998

$ws.Range("H121").Value = 2500
$ws.Range("J121").Value = 2500
$ws.Range("L121").Value = 7500
$ws.Range("N121").Value = -10994

$ws.Range("H129").Value = 2711.0435
$ws.Range("I129").Value = 1979.1666
$ws.Range("J129").Value = 3509.4546
$ws.Range("K129").Value = 5937.4998
$ws.Range("L129").Value = 10528.3638
$ws.Range("M129").Value = -937.4997999999996
$ws.Range("N129").Value = -20528.3638

$ws.Range("H137").Value = 2449.375
$ws.Range("I137").Value = 2157.9688
$ws.Range("J137").Value = 2837.9167
$ws.Range("K137").Value = 6473.9064
$ws.Range("L137").Value = 8513.750100000001
$ws.Range("M137").Value = -3923.9064
$ws.Range("N137").Value = -13613.7501

$ws.Range("H138").Value = 329350.72
$ws.Range("I138").Value = 3352.3845
$ws.Range("J138").Value = 395569.12
$ws.Range("K138").Value = 10057.1535
$ws.Range("L138").Value = 1186707.36
$ws.Range("M138").Value = -4917.1535
$ws.Range("N138").Value = -1196987.36

$ws.Range("H141").Value = 2013.4286
$ws.Range("I141").Value = 1219
$ws.Range("J141").Value = 3999.5
$ws.Range("K141").Value = 3657
$ws.Range("L141").Value = 11998.5
$ws.Range("M141").Value = 1523
$ws.Range("N141").Value = -22358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10621.153
$ws.Range("I32").Value = 6913.6353
$ws.Range("J32").Value = 34862.617
$ws.Range("K32").Value = 6913.6353
$ws.Range("L32").Value = 34862.617
$ws.Range("M32").Value = -6626.6353
$ws.Range("N32").Value = -35436.617

$ws.Range("H45").Value = 12200.533
$ws.Range("I45").Value = 15603.903
$ws.Range("K45").Value = 15603.903
$ws.Range("M45").Value = -15226.903

$ws.Range("H102").Value = 4556.231
$ws.Range("I102").Value = 4414.2354
$ws.Range("K102").Value = 4414.2354
$ws.Range("M102").Value = -2792.2354

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25255990
$ws.Range("I20").Value = 33337300
$ws.Range("K20").Value = 33337300
$ws.Range("M20").Value = -33337053

$ws.Range("H94").Value = 111111976
$ws.Range("I94").Value = 166667470
$ws.Range("K94").Value = 166667470
$ws.Range("M94").Value = -166667019

$ws.Range("H95").Value = 74979.25
$ws.Range("J95").Value = 74979.25
$ws.Range("L95").Value = 74979.25
$ws.Range("N95").Value = -80471.25

$ws.Range("H105").Value = 28892510
$ws.Range("I105").Value = 5002050
$ws.Range("K105").Value = 5002050
$ws.Range("M105").Value = -5000303

$ws.Range("H107").Value = 1397.5834
$ws.Range("I107").Value = 1243.625
$ws.Range("K107").Value = 1243.625
$ws.Range("M107").Value = 676.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6518.87
$ws.Range("I31").Value = 5196.387
$ws.Range("J31").Value = 11985.134
$ws.Range("K31").Value = 5196.387
$ws.Range("L31").Value = 11985.134
$ws.Range("M31").Value = -4901.387
$ws.Range("N31").Value = -12575.134

$ws.Range("H34").Value = 6518.87
$ws.Range("I34").Value = 5196.387
$ws.Range("J34").Value = 11985.134
$ws.Range("K34").Value = 5196.387
$ws.Range("L34").Value = 11985.134
$ws.Range("M34").Value = -4994.387
$ws.Range("N34").Value = -12389.134

$ws.Range("H75").Value = 49992.668
$ws.Range("J75").Value = 49992.668
$ws.Range("L75").Value = 49992.668
$ws.Range("N75").Value = -51988.668

$ws.Range("H78").Value = 49992.668
$ws.Range("J78").Value = 49992.668
$ws.Range("L78").Value = 149978.004
$ws.Range("N78").Value = -159962.004

$ws.Range("H119").Value = 761000
$ws.Range("J119").Value = 761000
$ws.Range("L119").Value = 761000
$ws.Range("N119").Value = -770676

$ws.Range("H120").Value = 37888.75
$ws.Range("J120").Value = 37888.75
$ws.Range("L120").Value = 37888.75
$ws.Range("N120").Value = -45146.75

$ws.Range("H132").Value = 11367390
$ws.Range("I132").Value = 12502928
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 37508784
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -37506254
$ws.Range("N132").Value = -41060

$ws.Range("H134").Value = 3613
$ws.Range("I134").Value = 2436.6086
$ws.Range("K134").Value = 7309.825800000001
$ws.Range("M134").Value = -4774.825800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 92.5
$ws.Range("I35").Value = 92.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 277.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 10.5
$ws.Range("N35").ClearContents()

$ws.Range("H69").Value = 1546
$ws.Range("J69").Value = 2368
$ws.Range("L69").Value = 7104
$ws.Range("N69").Value = -8726

$ws.Range("H72").Value = 1546
$ws.Range("J72").Value = 2368
$ws.Range("L72").Value = 21312
$ws.Range("N72").Value = -29424

$ws.Range("H86").Value = 378.1
$ws.Range("I86").Value = 296.66666
$ws.Range("J86").Value = 1111
$ws.Range("K86").Value = 889.9999799999999
$ws.Range("L86").Value = 3333
$ws.Range("M86").Value = 296.0000200000001
$ws.Range("N86").Value = -5705

$ws.Range("H89").Value = 378.1
$ws.Range("I89").Value = 296.66666
$ws.Range("J89").Value = 1111
$ws.Range("K89").Value = 2669.99994
$ws.Range("L89").Value = 9999
$ws.Range("M89").Value = 3258.00006
$ws.Range("N89").Value = -21855

$ws.Range("H113").Value = 6950
$ws.Range("J113").Value = 6950
$ws.Range("L113").Value = 20850
$ws.Range("N113").Value = -25190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22828364
$ws.Range("I70").Value = 29534828
$ws.Range("K70").Value = 29534828
$ws.Range("M70").Value = -29534558

$ws.Range("H73").Value = 22828364
$ws.Range("I73").Value = 29534828
$ws.Range("K73").Value = 29534828
$ws.Range("M73").Value = -29533892

$ws.Range("H80").Value = 333337730
$ws.Range("J80").Value = 6583
$ws.Range("L80").Value = 6583
$ws.Range("N80").Value = -8579

$ws.Range("H83").Value = 333337730
$ws.Range("J83").Value = 6583
$ws.Range("L83").Value = 32915
$ws.Range("N83").Value = -42899

$ws.Range("H97").Value = 718
$ws.Range("I97").Value = 756.2308
$ws.Range("K97").Value = 756.2308
$ws.Range("M97").Value = -260.2308

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1643.2667
$ws.Range("J46").Value = 1682.0714
$ws.Range("L46").Value = 1682.0714
$ws.Range("N46").Value = -2058.0714

$ws.Range("H122").Value = 7524.1177
$ws.Range("I122").Value = 7000
$ws.Range("J122").Value = 8113.75
$ws.Range("K122").Value = 21000
$ws.Range("L122").Value = 24341.25
$ws.Range("M122").Value = -18550
$ws.Range("N122").Value = -29241.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2444.67
$ws.Range("I132").Value = 2475.9524
$ws.Range("J132").Value = 2280.4375
$ws.Range("K132").Value = 7427.8572
$ws.Range("L132").Value = 6841.3125
$ws.Range("M132").Value = -4897.8572
$ws.Range("N132").Value = -11901.3125
